# Applies "Added most parts of cavetown" commit:
#  - Todo sheet: add a new TODO note about events 372/373
#  - GotoPoints sheet: add 11 new cave-town goto point descriptions
#  - Selection / active-sheet bookkeeping to match the authored session

$wb = $excel.ActiveWorkbook

# --- Todo sheet -----------------------------------------------------------
$todo = $wb.Worksheets.Item("Todo")
$todo.Range("A5").Value = "I guess I added wrong events to 372 or 373"
$todo.Range("B5").Select() | Out-Null

# --- GotoPoints sheet -------------------------------------------------------
$goto = $wb.Worksheets.Item("GotoPoints")

$newGotoPoints = @(
    "79: Warenhändler / Good merchant (Cavetown)",
    "80: Schmied / Blacksmith (Cavetown)",
    "81: Cavetown Büro / Cavetown Office (Cavetown)",
    "82: Vielauge-Schloss / Manyeyes' Castle (Cavetown)",
    "83: Gasthaus / Tavern (Cavetown)",
    "84: Badehaus / Bathhouse (Cavetown)",
    "85: Flosshändler / Raft Dealer (Cavetown)",
    "86: Gasthaus / Tavern (Cavetown) -- Second door",
    "87: Stadthaus 1 / Townhouse 1 (Cavetown)",
    "88: Stadthaus 2 / Townhouse 2 (Cavetown)",
    "89: Stadthaus 3 / Townhouse 3 (Cavetown)"
)

$row = 6
foreach ($text in $newGotoPoints) {
    $goto.Cells.Item($row, 1).Value = $text
    $row = $row + 1
}

$goto.Activate() | Out-Null
$goto.Range("C15").Select() | Out-Null
